$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192046999931335
$ws.Range("B1").Value = 2.457992792129517
$ws.Range("C1").Value = 3.952888965606689
$ws.Range("D1").Value = 2.143872737884521
$ws.Range("E1").Value = 1.189554333686829
